$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B24").Value = "LOM3229 -  Métodos Experimentais da Física II  (Indicação de Conjunto)`n"
$ws.Range("C24").Value = "LOM3229 -  Métodos Experimentais da Física II  (Indicação de Conjunto)`n"

$ws.Range("B25").Value = "LOB1021 -  Física IV  (Requisito)`n"
$ws.Range("C25").Value = "LOB1021 -  Física IV  (Requisito)`n"

$ws.Range("B26").Value = "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)`n"
$ws.Range("C26").Value = "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)`n"
